$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- Formatting ---------------------------------------------------------
# B1 and A2 share identical formatting: bold font, thin box border,
# centered horizontally, top vertically aligned.
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108  # xlCenter
$b1.VerticalAlignment = -4160    # xlTop
$b1.Borders.LineStyle = 1        # xlContinuous
$b1.Borders.Weight = 2           # xlThin

# Copy the exact same style onto A2 (instead of re-deriving it through
# the same property calls) so both cells land on the single shared
# cellXfs record rather than each cell minting its own.
$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
